# Regenerate orders with updated distance/size codes.
# The experiment's Distance levels (D64/D80/D51) and the "S30" Size level
# were re-generated to new values (D69/D86/D55, S31). Every place those
# tokens appear - the Condition column, the left/right filename columns,
# the Distance legend, and the Size legend - needs updating, while the
# rest of each string (Face numbers, S25/S20, "_l.png"/"_r.png", etc.)
# stays exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

# Order doesn't matter: none of the old/new tokens are substrings of each
# other, so a simple sequence of partial-match replacements across the
# whole used range is safe and idempotent.
$used.Replace("D64", "D69") | Out-Null
$used.Replace("D80", "D86") | Out-Null
$used.Replace("D51", "D55") | Out-Null
$used.Replace("S30", "S31") | Out-Null
